# Change in unit of AIC and run of new results.
#
# Every yearly worksheet ("2000" .. "2100") shares the same small table:
# rows 5, 7 and 8 hold numeric results in columns D, E, F and G.
# The underlying model was re-run with the AIC figures expressed in a
# different (smaller) unit, which is equivalent to scaling every one of
# those cached numbers by 1e-6 (e.g. mg -> kg). Apply that conversion to
# every sheet in the workbook.

$wb = $excel.ActiveWorkbook
$factor = [double]"0.000001"
$cols = @("D", "E", "F", "G")
$rows = @(5, 7, 8)

foreach ($ws in $wb.Worksheets) {
    foreach ($r in $rows) {
        foreach ($c in $cols) {
            $cell = $ws.Range("$c$r")
            $old = $cell.Value2
            $cell.Value = $old * $factor
        }
    }
}
